$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.104.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.98%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.638.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.03%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'216.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.03%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.518"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.91%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.14%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.17%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.55%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.11%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.867.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.05%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.638.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.12%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +0.05%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.542"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.38%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.73%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'27.114.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.06%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.0₃0740"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.34%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.26%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'6.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.69%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.57%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +3.21%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.59%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'146.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.38%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.03%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +1.39%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.32%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.33%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0507"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.08%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.15%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.63%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.73%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.308.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.08%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +1.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.49%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +2.90%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.542"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.73%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.00%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.50%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +5.77%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.45%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.777.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.05%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'61.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.04%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'91.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.29%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.82%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0₆0108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.19%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.09%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'7.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.17%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0960"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.22%  "
$ws.Range("E51").Style = "Normal"
